# Insert a new row at position 16 (pushes existing rows 16-36 down to 17-37)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

# Fill in the new row 16 with the new weekly record
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44477
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112022
$ws.Range("G16").Value = "Arveja Verde"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 23000
$ws.Range("L16").Value = 23000
$ws.Range("M16").Value = 23000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región de O'Higgins"
$ws.Range("P16").Value = 920
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"

# Make sure the date cell keeps the date number format (style index 2 from styles.xml)
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
